# edit.ps1
# Applies the "DEXAH_output" edit: updates text on existing slides 2-6,
# converts slide 5 from a title-only layout to a title-slide (ctrTitle/subTitle)
# layout, and appends 23 new slides (7-29) with lyrics content.

$p = $ppt.ActivePresentation
$NL = [char]13

# ---------------------------------------------------------------------------
# Slide 2: "CON CANTICOS, SENOR" title slide -> "TE AMO, REY"
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "TE AMO, REY"
$s2.Shapes.Item(2).TextFrame.TextRange.Text = "46 - Himnario Majestuoso"

# ---------------------------------------------------------------------------
# Slide 3
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Te amo Rey, y levanto mi voz, para adorarte, mi salvador. Me gozo en ti y te alabo mi Dios; dulce sea mi canto a ti, oh Señor."

# ---------------------------------------------------------------------------
# Slide 4
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Te amo Rey, y levanto mi voz, para adorar y gozar me en ti. Regocíjate y escucha, mi Rey: que sea un dulce sonar para ti."

# ---------------------------------------------------------------------------
# Slide 5: was a "titleOnly" lyrics slide; becomes a title slide
# (ctrTitle/subTitle) like slide 2. The runtime's shape Delete()/Layout
# setter do not remove the old placeholder, so recreate the slide by
# deleting it and inserting a fresh one with the title layout in its place.
# ---------------------------------------------------------------------------
$p.Slides.Item(5).Delete()
$s5 = $p.Slides.Add(5, 1)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Dios descendió"
$s5.Shapes.Item(2).TextFrame.TextRange.Text = "333 - Himnario Majestuoso"

# ---------------------------------------------------------------------------
# Slide 6
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Día tan grande no puedo olvidar, día de gloria sin par; cuando en tinieblas al verme andar, vino a salvar me el señor." + $NL + " Gran compasión tuvo Cristo de mí, de gozo y paz me llenó. Quitó las sombras, oh gloria a su nombre, la noche en día a cambió."

# ---------------------------------------------------------------------------
# New slides 7-29 (titleOnly = layout 11, title/ctrTitle+subTitle = layout 1)
# ---------------------------------------------------------------------------

# Slide 7
$s = $p.Slides.Add(7, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Dios descendió, de gloria me llenó, (me llenó) cuando Jesús por gracia me salvó; (me salvó) fui ciego, me hizo ver, y en él renacer. Dios descendió y de gloria me llenó. (me llenó)"

# Slide 8
$s = $p.Slides.Add(8, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Nací de nuevo en virtud de Jesús, a la familia de Dios. Justificado por Cristo el señor, gozo la gran redención." + $NL + " Bendito sea mi padre y Dios, que cuando vine con fe, fui adoptado por Cristo el amado, loores por siempre daré."

# Slide 9
$s = $p.Slides.Add(9, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Dios descendió, de gloria me llenó, (me llenó) cuando Jesús por gracia me salvó; (me salvó) fui ciego, me hizo ver, y en él renacer. Dios descendió y de gloria me llenó. (me llenó)"

# Slide 10 (title slide)
$s = $p.Slides.Add(10, 1)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Fue Tu gracia"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "48 - Himnario Gracia"

# Slide 11
$s = $p.Slides.Add(11, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "  Siendo yo un pecador, Extraviado en el error; Incapaz fui de quitar la culpa que cargaba Mi corazón. "

# Slide 12
$s = $p.Slides.Add(12, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = " Fue Tu gracia Que mi vida renovó, Fue Tu gracia Que mi deuda canceló, Fue Tu gracia Que de muerte me llevó A vivir para Tu gloria, oh Señor. "

# Slide 13
$s = $p.Slides.Add(13, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = " ¿Tu gracia, quién puede medir? Tu amor no tiene uno igual. Cuando perdido estaba yo viniste a mi encuentro, Mi Salvador.(Coro) "

# Slide 14
$s = $p.Slides.Add(14, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = " Fue Tu gracia Que mi vida renovó, Fue Tu gracia Que mi deuda canceló, Fue Tu gracia Que de muerte me llevó A vivir para Tu gloria, oh Señor. "

# Slide 15
$s = $p.Slides.Add(15, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = " En la cruz donde murió El perfecto Hijo de Dios, La sentencia Él sufrió, la que yo merecía Por mi maldad. (Coro) "

# Slide 16
$s = $p.Slides.Add(16, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = " Fue Tu gracia Que mi vida renovó, Fue Tu gracia Que mi deuda canceló, Fue Tu gracia Que de muerte me llevó A vivir para Tu gloria, oh Señor. "

# Slide 17
$s = $p.Slides.Add(17, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = " Tu gracia me ha traído aquí, De rodillas ante Ti; Maravillado al comprender que a Tu Hijo entregaste En mi lugar. (Coro)  "

# Slide 18
$s = $p.Slides.Add(18, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = " Fue Tu gracia Que mi vida renovó, Fue Tu gracia Que mi deuda canceló, Fue Tu gracia Que de muerte me llevó A vivir para Tu gloria, oh Señor. "

# Slide 19 (empty paragraph, leave placeholder text untouched)
$s = $p.Slides.Add(19, 11)

# Slide 20
$s = $p.Slides.Add(20, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = " Fue Tu gracia Que mi vida renovó, Fue Tu gracia Que mi deuda canceló, Fue Tu gracia Que de muerte me llevó A vivir para Tu gloria, oh Señor. "

# Slide 21 (title slide)
$s = $p.Slides.Add(21, 1)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Mi vida es Cristo"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "143 - Himnario Gracia"

# Slide 22
$s = $p.Slides.Add(22, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "  En densa oscuridad vagué, perdido en el error; La senda vana del placer a muerte me llevó. Siendo rebelde a Tu voz quisiste amarme así; De no haber sido por Tu amor aún huiría de Ti. "

# Slide 23
$s = $p.Slides.Add(23, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = " //¡Aleluya! Mi vida es Cristo. ¡Aleluya! Jesús es mi todo.//  "

# Slide 24
$s = $p.Slides.Add(24, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = " En rumbo a mi perdición, indiferente aún. De mí tuviste compasión, me guiaste a la cruz; Y contemplé Tu gran bondad: sufriste Tú por mí; Al Tú morir en mi lugar Tu gracia recibí. "

# Slide 25
$s = $p.Slides.Add(25, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = " //¡Aleluya! Mi vida es Cristo. ¡Aleluya! Jesús es mi todo.//  "

# Slide 26
$s = $p.Slides.Add(26, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = " Ahora Señor Tuyo seré, y viviré por Ti. Tus mandamientos seguiré por Tu poder en mí. Usa mi vida, oh Señor, como lo quieras Tú; Y que sea siempre mi canción: ""Mi gloria eres Tú."" "

# Slide 27
$s = $p.Slides.Add(27, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = " //¡Aleluya! Mi vida es Cristo. ¡Aleluya! Jesús es mi todo.//  "

# Slide 28 (empty paragraph, leave placeholder text untouched)
$s = $p.Slides.Add(28, 11)

# Slide 29
$s = $p.Slides.Add(29, 11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = " //¡Aleluya! Mi vida es Cristo. ¡Aleluya! Jesús es mi todo.//  "
